# Backup QR Scanner data - 17/05/2025, 10:46:23 PM
# Append one new scan record as row 2 of the Checklist sheet (A1:F1 -> A1:F2),
# keeping every field as literal text (matching the existing header row's
# text-only storage) rather than letting Excel auto-coerce number/date-like
# values into numeric/date cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = @{
    A = "228723"
    B = "Biochemistry"
    C = "05/17/2025"
    D = "22:44:27"
    E = "Manual"
    F = "231249@med.asu.edu.eg"
}

$rowIndex = 2
foreach ($col in "A", "B", "C", "D", "E", "F") {
    $value = $newRow[$col]
    $cell = $ws.Range("$col$rowIndex")

    # Values that look like a number or a date/time get auto-converted by
    # Excel on plain assignment (e.g. "228723" -> 228723, "05/17/2025" ->
    # a date serial). Prefix those with an apostrophe so they stay text,
    # exactly like the other literal-text values already on the sheet.
    if ($value -match '^[0-9]+(\.[0-9]+)?$' -or $value -match '^\d{1,2}/\d{1,2}/\d{2,4}$' -or $value -match '^\d{1,2}:\d{2}(:\d{2})?$') {
        $cell.Value = "'" + $value
    } else {
        $cell.Value = $value
    }
}

# The apostrophe prefix also stamps the cell with a "quote prefix" style so
# Excel can redraw the little text-alignment indicator; that leaves a stray
# style index that the header row doesn't have. Re-apply the header row's
# (unstyled) formatting onto the new row so row 2 matches row 1's plain
# formatting, without touching the values/text we just wrote.
$ws.Range("A1:F1").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
